# Create Entity & Send Notification
# Updates the "tracking number" / date / confirmation-timestamp values that
# are stamped onto row 2 of the first four worksheets (Sheet1..Sheet4) each
# time the "Create Entity" + "Send Notification" scenario is re-run.
#
# Columns G, AF, AV, AZ hold an SMS/tracking reference number (digits only -
# must stay TEXT so Excel doesn't coerce it to a numeric value and lose any
# formatting). Columns O and AD hold a plain date string, which Excel would
# otherwise auto-detect and convert to a date serial, so those are also
# forced to Text before being written. Columns Q and BB hold strings that
# Excel's literal-entry parser does not recognise as a date/number (they
# contain words / AM-PM text), so they can be written directly.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Range, [string]$Value)
    # Force Text number format first so digit-only / date-look-alike
    # strings are kept as literal text instead of being coerced into a
    # number or date serial value by Excel's input parser.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
}

$trackingNumber = "9840038822"
$trackingNumberAF = "9840011089"
$trackingNumberAV = "9840036820"
$trackingNumberAZ = "9840019753"
$dateValue = "06-01-2025"
$dueDateValue = "09-01-2025 05:00:00 PM"

# Sheet1 ("Create Entity" row)
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1.Range("G2") $trackingNumber
Set-TextValue $ws1.Range("O2") $dateValue
$ws1.Range("Q2").Value = $dueDateValue
Set-TextValue $ws1.Range("AD2") $dateValue
Set-TextValue $ws1.Range("AF2") $trackingNumberAF
Set-TextValue $ws1.Range("AV2") $trackingNumberAV
Set-TextValue $ws1.Range("AZ2") $trackingNumberAZ
$ws1.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 1:12 PM"

# Sheet2 ("Summary" row)
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("G2") $trackingNumber
Set-TextValue $ws2.Range("O2") $dateValue
$ws2.Range("Q2").Value = $dueDateValue
Set-TextValue $ws2.Range("AD2") $dateValue
Set-TextValue $ws2.Range("AF2") $trackingNumberAF
Set-TextValue $ws2.Range("AV2") $trackingNumberAV
Set-TextValue $ws2.Range("AZ2") $trackingNumberAZ
$ws2.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 1:24 PM"

# Sheet3 ("Duplicate Email" row)
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("G2") $trackingNumber
Set-TextValue $ws3.Range("AF2") $trackingNumberAF
Set-TextValue $ws3.Range("AV2") $trackingNumberAV
Set-TextValue $ws3.Range("AZ2") $trackingNumberAZ

# Sheet4 ("Edit New Email" row)
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("G2") $trackingNumber
Set-TextValue $ws4.Range("AF2") $trackingNumberAF
Set-TextValue $ws4.Range("AV2") $trackingNumberAV
Set-TextValue $ws4.Range("AZ2") $trackingNumberAZ
